$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear rows that are fully removed or being rebuilt (21-22 old LPN-failed/terminated rows, 25 old NOTE row)
$ws.Range("A21:H22").Clear()
$ws.Range("A25").Clear()

# Rewrite the table contents (row 18 is newly inserted; rows below shift up by one)
$ws.Range("A1").Value = 'Server Model (Low Power Node)'
$ws.Range("A3").Value = 'Event Handler Event'
$ws.Range("B3").Value = 'Condition'
$ws.Range("C3").Value = 'Command(s)'
$ws.Range("D3").Value = 'Resulting Events'
$ws.Range("E3").Value = 'DISPLAY_ROW_ACTION update'
$ws.Range("F3").Value = 'DISPLAY_ROW_LPN update'
$ws.Range("G3").Value = 'DISPLAY_ROW_CONNECTION update'
$ws.Range("H3").Value = 'Notes'
$ws.Range("A4").Value = 'gecko_evt_system_boot_id'
$ws.Range("B4").Value = 'button 0 or button 1 pressed'
$ws.Range("C4").Value = 'gecko_cmd_flash_ps_erase_all; gecko_cmd_hardware_set_soft_timer'
$ws.Range("D4").Value = 'gecko_evt_hardware_soft_timer_id'
$ws.Range("E4").Value = 'Factory Reset'
$ws.Range("H4").Value = 'Factory reset and clear provisioner data for test purposes when either button is held on startup.  Must configure a timer to perform reset after 2 seconds'
$ws.Range("A5").Value = 'gecko_evt_system_boot_id'
$ws.Range("B5").Value = 'neither button 0 or button 1 pressed'
$ws.Range("C5").Value = 'SetDeviceName; gecko_cmd_mesh_node_init; gecko_cmd_flash_ps_load'
$ws.Range("D5").Value = 'gecko_evt_mesh_node_initialized'
$ws.Range("H5").Value = 'Initialize the mesh stack'
$ws.Range("A6").Value = 'gecko_evt_hardware_soft_timer_id'
$ws.Range("B6").Value = 'factory reset timer handle, display update timer handle, log time timer handle, friend find timer handle, alerts handle, interrupt drive sensor timeout events'
$ws.Range("C6").Value = 'gecko_cmd_system_reset, displayUpdate, gecko_cmd_mesh_lpn_establish_friendship, gecko_cmd_flash_ps_save'
$ws.Range("D6").Value = 'gecko_evt_mesh_lpn_friendship_established_id, gecko_evt_mesh_lpn_friendship_failed_id'
$ws.Range("H6").Value = 'Completes factory reset, updates display for remove charge buildup, log timestamp value increase, tries to establish friendship on Client side'
$ws.Range("A7").Value = 'gecko_evt_mesh_node_initialized_id'
$ws.Range("B7").Value = 'provisioned'
$ws.Range("C7").Value = 'gecko_cmd_mesh_generic_server_init, gecko_cmd_mesh_friend_init(), gpioIntEnable(), mesh_lib_init(malloc,free,9), mesh_lib_generic_server_register_handler'
$ws.Range("D7").Value = 'N/A'
$ws.Range("H7").Value = 'Initialize the server model and friend initialization; also takes action based on saved persistent data'
$ws.Range("A8").Value = 'gecko_evt_mesh_node_initialized_id'
$ws.Range("B8").Value = '!provisioned'
$ws.Range("C8").Value = 'gecko_cmd_mesh_node_start_unprov_beaconing(0x03)'
$ws.Range("D8").Value = 'gecko_evt_mesh_node_provisioning_started, gecko_evt_mesh_node_provisioned, gecko_evt_mesh_node_provisioning_failed'
$ws.Range("H8").Value = 'Start beaconing for provisioning in PB-ADV and PB-GATT modes'
$ws.Range("A9").Value = 'gecko_evt_mesh_node_provisioning_started_id'
$ws.Range("B9").Value = 'N/A'
$ws.Range("C9").Value = 'N/A'
$ws.Range("D9").Value = 'N/A'
$ws.Range("E9").Value = 'Provisioning'
$ws.Range("A10").Value = 'gecko_evt_mesh_node_provisioned_id'
$ws.Range("B10").Value = 'N/A'
$ws.Range("C10").Value = 'N/A'
$ws.Range("D10").Value = 'N/A'
$ws.Range("E10").Value = 'Provisioned'
$ws.Range("A11").Value = 'gecko_evt_mesh_node_provisioning_failed_id'
$ws.Range("B11").Value = 'N/A'
$ws.Range("C11").Value = 'N/A'
$ws.Range("D11").Value = 'N/A'
$ws.Range("E11").Value = 'Provision Fail'
$ws.Range("A12").Value = 'gecko_evt_mesh_generic_server_client_request_id'
$ws.Range("B12").Value = 'mesh_generic_request is received'
$ws.Range("C12").Value = 'mesh_lib_generic_server_event_handler'
$ws.Range("D12").Value = 'N/A'
$ws.Range("H12").Value = 'Deseralizes data when mesh generic request is received from friend'
$ws.Range("A13").Value = 'gecko_evt_mesh_generic_server_state_changed_id'
$ws.Range("C13").Value = 'mesh_lib_generic_server_event_handler'
$ws.Range("D13").Value = 'N/A'
$ws.Range("A14").Value = 'gecko_evt_le_connection_opened_id'
$ws.Range("B14").Value = 'N/A'
$ws.Range("C14").Value = 'N/A'
$ws.Range("D14").Value = 'N/A'
$ws.Range("G14").Value = 'Connected'
$ws.Range("A15").Value = 'gecko_evt_le_connection_closed_id'
$ws.Range("B15").Value = 'N/A'
$ws.Range("C15").Value = 'N/A'
$ws.Range("D15").Value = 'N/A'
$ws.Range("G15").Value = '  (blank)'
$ws.Range("A16").Value = 'gecko_evt_mesh_node_reset_id'
$ws.Range("B16").Value = 'gecko_cmd_hardware_set_soft_timer'
$ws.Range("D16").Value = 'gecko_evt_hardware_soft_timer_id'
$ws.Range("H16").Value = 'Factory reset when requested by provisioner, using the same 2 second soft timer delay'
$ws.Range("A17").Value = 'gecko_evt_system_external_signal_id'
$ws.Range("B17").Value = 'Flame sensor  interrupt and Gas sensor scheduler events and push button interrupt'
$ws.Range("C17").Value = 'mesh_lib_generic_server_publish, gecko_cmd_flash_ps_save, gecko_cmd_hardware_set_soft_timer'
$ws.Range("D17").Value = 'gecko_evt_hardware_soft_timer_id'
$ws.Range("H17").Value = 'Handle gas sensor scheduler events; fire sensor interrupt events; push button 0 interrupt and publishes required data; Also saves in persistent data'
$ws.Range("A18").Value = 'gecko_evt_mesh_lpn_friendship_established_id'
$ws.Range("B18").Value = 'N/A'
$ws.Range("C18").Value = 'N/A'
$ws.Range("D18").Value = 'N/A'
$ws.Range("F18").Value = 'LPN'
$ws.Range("H18").Value = 'Friendship is established when this node'
$ws.Range("A19").Value = 'gecko_evt_mesh_lpn_friendship_failed_id'
$ws.Range("B19").Value = 'gecko_cmd_mesh_lpn_establish_friendship'
$ws.Range("C19").Value = 'gecko_cmd_hardware_set_soft_timer'
$ws.Range("D19").Value = 'gecko_evt_hardware_soft_timer_id'
$ws.Range("F19").Value = 'no friend'
$ws.Range("H19").Value = 'Tries to establish friendship every 2 seconds until it is established'
$ws.Range("A20").Value = 'gecko_evt_mesh_lpn_friendship_terminated_id'
$ws.Range("B20").Value = 'gecko_cmd_mesh_lpn_establish_friendship'
$ws.Range("C20").Value = 'gecko_cmd_hardware_set_soft_timer'
$ws.Range("D20").Value = 'gecko_evt_hardware_soft_timer_id'
$ws.Range("F20").Value = 'friend lost'
$ws.Range("H20").Value = 'Tries to establish friendship every 2 seconds until it is established'

# A1 title is bold; A23 is a new empty bold placeholder cell (formerly the NOTE cell at A25)
$ws.Range("A1").Font.Bold = $true
$ws.Range("A23").Font.Bold = $true

# Restore the active selection shown in the saved file
$ws.Range("B13").Select()
